$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Form Tag" column is appended after the existing "Client Email"
# column (A:I), defaulting every existing investor row's tag to "Default".
$ws.Range("J1").Value = "Form Tag"

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 10).Value = "Default"
}

$ws.Range("J1").Select()
